$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")

# Mark "Día 6" (column W) as completed for the CU 22 (row 23) and CU 08 (row 24)
# task rows - this ripples the shared "remaining" formulas across the rest of
# the timeline (columns X..AY) and the totals in AZ/BA back down to 0 / up to 1.
$ws.Range("W23").Value = 1
$ws.Range("W24").Value = 1

# Update the view: the sheet had scrolled/selected to reflect where the user
# was working.
$ws.Activate()
$ws.Range("W23").Select()
